$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (victims): update existing rows 2-4, add row 5 ---
$ws1.Range("D2").Value = 6.049624213818273
$ws1.Range("E2").Value = 3.173920439702034
$ws1.Range("F2").Value = -0.04606928475336147
$ws1.Range("G2").Value = 0.01562672959107125
$ws1.Range("H2").Value = 0.04864744263986305
$ws1.Range("I2").Value = 57

$ws1.Range("D3").Value = 1.439783477814177
$ws1.Range("E3").Value = 1.251316609725247
$ws1.Range("F3").Value = 0.1119952863246196
$ws1.Range("G3").Value = 0.1190955254970143
$ws1.Range("H3").Value = 0.163482991018465
$ws1.Range("I3").Value = 0

$ws1.Range("D4").Value = 7.474574587849915
$ws1.Range("E4").Value = 8.261890110033912
$ws1.Range("F4").Value = -0.03158732575226431
$ws1.Range("G4").Value = -0.0401545081878627
$ws1.Range("H4").Value = 0.05108956523585606
$ws1.Range("I4").Value = 72

$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = 5.853086206137439
$ws1.Range("C5").Value = 6.835972487871987
$ws1.Range("D5").Value = 5.847429627558846
$ws1.Range("E5").Value = 6.824481192067192
$ws1.Range("F5").Value = -0.005656578578593141
$ws1.Range("G5").Value = -0.01149129580479524
$ws1.Range("H5").Value = 0.01280807403511955
$ws1.Range("I5").Value = 49

# --- Sheet2 (false_positives): update rows 2-11, add rows 12-21 ---
$ws2.Range("B2").Value = 5.348999931723383
$ws2.Range("C2").Value = 8.480579390302147

$ws2.Range("B3").Value = 7.526828432972257
$ws2.Range("C3").Value = 1.021908001361185

$ws2.Range("B4").Value = 6.837243571439553
$ws2.Range("C4").Value = 2.405244964820472

$ws2.Range("B5").Value = 7.905431378799093
$ws2.Range("C5").Value = 5.331689761992734

$ws2.Range("B6").Value = 3.397695124299078
$ws2.Range("C6").Value = 4.381497769581268

$ws2.Range("B7").Value = 1.226557369163704
$ws2.Range("C7").Value = 1.994266211996512

$ws2.Range("B8").Value = 6.364995317549043
$ws2.Range("C8").Value = 6.177516092594001

$ws2.Range("B9").Value = 5.923080891850031
$ws2.Range("C9").Value = 4.069420434095067

$ws2.Range("B10").Value = 8.977679486313688
$ws2.Range("C10").Value = 8.846682710209841

$ws2.Range("B11").Value = 2.08077204017929
$ws2.Range("C11").Value = 6.771906721552654

$ws2.Range("A12").Value = 11
$ws2.Range("B12").Value = 5.202834579805807
$ws2.Range("C12").Value = 3.481935004471645

$ws2.Range("A13").Value = 12
$ws2.Range("B13").Value = 8.472348127649997
$ws2.Range("C13").Value = 3.862361573672562

$ws2.Range("A14").Value = 13
$ws2.Range("B14").Value = 4.13295200422529
$ws2.Range("C14").Value = 8.122194816038338

$ws2.Range("A15").Value = 14
$ws2.Range("B15").Value = 2.817260748267038
$ws2.Range("C15").Value = 5.985497157488339

$ws2.Range("A16").Value = 15
$ws2.Range("B16").Value = 1.672122748659079
$ws2.Range("C16").Value = 7.661153181227182

$ws2.Range("A17").Value = 16
$ws2.Range("B17").Value = 3.688936484365283
$ws2.Range("C17").Value = 2.202235735158713

$ws2.Range("A18").Value = 17
$ws2.Range("B18").Value = 4.602714933194296
$ws2.Range("C18").Value = 7.370594162298354

$ws2.Range("A19").Value = 18
$ws2.Range("B19").Value = 2.845137671949979
$ws2.Range("C19").Value = 1.416170408515277

$ws2.Range("A20").Value = 19
$ws2.Range("B20").Value = 1.726024364952975
$ws2.Range("C20").Value = 5.642659087894805

$ws2.Range("A21").Value = 20
$ws2.Range("B21").Value = 2.596123551745706
$ws2.Range("C21").Value = 8.536904884051982
